# Applies the "Microsite Education Script completed" update:
#  - AMSIN  sheet: append rows 53-57 (165_fstcycle ... 166cyclescnd)
#  - BETA   sheet: append rows 26-27 (165beta, 166_beta)
#  - AMS    sheet: fix row 25 (B value correction) and append rows 26-27
#                  (165_live, 166_live)

$wb = $excel.ActiveWorkbook

# A reusable scratch cell (formatted as Text) used as a "format donor" so that
# date-looking strings such as "2022-08-02" get written as literal text
# (matching the source inlineStr cells) instead of being auto-converted by
# Excel's smart entry into a date serial number.
function Get-TextScratch($ws) {
    $scratch = $ws.Cells.Item(5000, 26)   # far-away unused cell (column Z)
    $scratch.NumberFormat = "@"
    return $scratch
}

function Set-TextValue($ws, $scratch, $row, $col, $text) {
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-DateValue($ws, $donorCell, $row, $col, $serial) {
    # Copies the number-format (and only the number-format) from a cell that
    # already carries the "YYYY-MM-DD HH:MM:SS" style used throughout column B
    $donorCell.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, $col).Value = $serial
}

function Set-DataRow($ws, $scratch, $donorCell, $row, $a, $b, $c, $d, $e, $f, $g) {
    Set-TextValue $ws $scratch $row 1 $a
    Set-DateValue $ws $donorCell $row 2 $b
    Set-TextValue $ws $scratch $row 3 $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
}

# ---------------------------------------------------------------------------
# Sheet "AMSIN" -> add rows 53-57
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
$scratchAmsin = Get-TextScratch $wsAmsin
$donorAmsin = $wsAmsin.Range("B52")   # existing cell using the target date style (s=12)

Set-DataRow $wsAmsin $scratchAmsin $donorAmsin 53 "2022-08-02" 44775.66379868056 "165_fstcycle" 96 95 1 2.6
Set-DataRow $wsAmsin $scratchAmsin $donorAmsin 54 "2022-08-03" 44776.69014003472 "165_scndcycle" 96 95 1 2.89
Set-DataRow $wsAmsin $scratchAmsin $donorAmsin 55 "2022-08-04" 44777.39273803241 "165_finalrun" 96 95 1 2.58
Set-DataRow $wsAmsin $scratchAmsin $donorAmsin 56 "2022-08-22" 44795.67178451389 "166fstcycle" 96 92 4 3.27
Set-DataRow $wsAmsin $scratchAmsin $donorAmsin 57 "2022-08-23" 44796.90697054398 "166cyclescnd" 96 96 0 2.49

$scratchAmsin.Clear() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "BETA" -> add rows 26-27
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
$scratchBeta = Get-TextScratch $wsBeta
$donorBeta = $wsBeta.Range("B25")   # existing cell using the target date style (s=12)

Set-DataRow $wsBeta $scratchBeta $donorBeta 26 "2022-08-04" 44777.56231642361 "165beta" 96 96 0 2.44
Set-DataRow $wsBeta $scratchBeta $donorBeta 27 "2022-08-24" 44797.54109481481 "166_beta" 96 96 0 2.4

$scratchBeta.Clear() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "AMS" -> correct row 25 (B value) and add rows 26-27
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")
$scratchAms = Get-TextScratch $wsAms
$donorAms = $wsAms.Range("B24")   # existing cell using the target date style (s=11)

# Row 25 already holds the correct text/numbers; only the timestamp in B25
# needs correcting, and the run re-applies the same look (general style)
# consistently across the row.
Set-DataRow $wsAms $scratchAms $donorAms 25 "2022-07-14" 44756.82067131944 "164_live" 96 96 0 2.38

Set-DataRow $wsAms $scratchAms $donorAms 26 "2022-08-04" 44777.81690030092 "165_live" 96 86 10 4.13
Set-DataRow $wsAms $scratchAms $donorAms 27 "2022-08-24" 44797.92316747576 "166_live" 96 96 0 2.5

$scratchAms.Clear() | Out-Null
